# Refresh the cryptocurrency price / 1h-volume snapshot (scheduled GitHub Actions run).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# cell reference -> new text. Values are written through NumberFormat "@" (Text)
# + a Style reset to "Normal" so numeric-looking strings such as "213.15", "1.00" or
# "0.0₇0956" are kept as literal text (matching the source inlineStr cells) instead
# of being reinterpreted as numbers, while leaving no residual style on the cell.
$updates = @(
    @("D2", "26.385.45"),
    @("E2", "  +1.03%  "),
    @("D3", "1.611.45"),
    @("E3", "  +0.74%  "),
    @("E4", "  -0.37%  "),
    @("D5", "213.15"),
    @("E5", "  +0.31%  "),
    @("E6", "  -0.51%  "),
    @("E7", "  +0.68%  "),
    @("E8", "  +0.74%  "),
    @("D9", "0.0615"),
    @("E9", "  +0.09%  "),
    @("D10", "18.47"),
    @("E10", "  +2.90%  "),
    @("D11", "0.0813"),
    @("E11", "  -0.82%  "),
    @("D12", "1.836.42"),
    @("E12", "  +0.73%  "),
    @("D13", "1.621.49"),
    @("E13", "  +1.24%  "),
    @("D14", "4.04"),
    @("E14", "  +1.09%  "),
    @("E15", "  +1.31%  "),
    @("D16", "26.367.73"),
    @("E16", "  +1.00%  "),
    @("D17", "62.15"),
    @("E17", "  +2.60%  "),
    @("E18", "  +1.42%  "),
    @("E19", "  -0.30%  "),
    @("D20", "203.28"),
    @("E20", "  -2.20%  "),
    @("E21", "  +1.58%  "),
    @("E22", "  +0.47%  "),
    @("D23", "6.03"),
    @("E23", "  +0.90%  "),
    @("E24", "  +2.70%  "),
    @("D25", "144.88"),
    @("E25", "  +2.26%  "),
    @("E26", "  -0.26%  "),
    @("D27", "0.122"),
    @("E27", "  -2.21%  "),
    @("D28", "15.23"),
    @("E28", "  -0.01%  "),
    @("E29", "  +2.35%  "),
    @("D30", "0.0492"),
    @("E30", "  +4.81%  "),
    @("E31", "  +0.52%  "),
    @("D32", "3.21"),
    @("E32", "  +2.57%  "),
    @("E33", "  -1.75%  "),
    @("E34", "  +2.66%  "),
    @("E35", "  +1.70%  "),
    @("D36", "1.168.06"),
    @("E36", "  +5.54%  "),
    @("E37", "  +3.15%  "),
    @("E38", "  -0.38%  "),
    @("D39", "0.790"),
    @("E39", "  +1.65%  "),
    @("E40", "  -0.10%  "),
    @("E41", "  +1.68%  "),
    @("D42", "0.783"),
    @("E42", "  +1.01%  "),
    @("E43", "  +3.49%  "),
    @("D44", "1.751.14"),
    @("E44", "  +0.52%  "),
    @("D45", "91.99"),
    @("E45", "  -0.66%  "),
    @("E46", "  +1.72%  "),
    @("D47", "54.35"),
    @("E47", "  +1.84%  "),
    @("D48", "0.0507"),
    @("E48", "  +0.36%  "),
    @("D49", "0.408"),
    @("E49", "  -0.55%  "),
    @("B50", "BabyDogeCoin"),
    @("C50", "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"),
    @("D50", "0.0₇0956"),
    @("E50", "  -10.49%  "),
    @("B51", "USDD"),
    @("C51", "https://coinranking.com/coin/z2PZIKQL7+usdd-usdd"),
    @("D51", "1.00"),
    @("E51", "  -0.30%  ")
)

foreach ($u in $updates) {
    $cellRef = $u[0]
    $newText = $u[1]
    $r = $ws.Range($cellRef)
    $r.NumberFormat = "@"
    $r.Value = $newText
    $r.Style = "Normal"
}
